$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2600345611572266
$ws.Range("E2").Value = 92.32123641794169
$ws.Range("F2").Value = 0.003757337172971522
$ws.Range("G2").Value = 0.002938912155312964
$ws.Range("H2").Value = 0.002563019735842009
$ws.Range("I2").Value = 0.002556876524694315
$ws.Range("J2").Value = 0.002467525098318024
$ws.Range("K2").Value = 0.002429754410296352
$ws.Range("L2").Value = 0.002193355671268605
$ws.Range("M2").Value = 0.002166346550085777
$ws.Range("N2").Value = 0.002075890163776998
$ws.Range("O2").Value = 0.002075890163776998
$ws.Range("P2").Value = 0.00204430077160284
$ws.Range("Q2").Value = 0.00204430077160284
$ws.Range("R2").Value = 0.001879124199714181
$ws.Range("S2").Value = 0.001879124199714181
$ws.Range("T2").Value = 0.001879124199714181
$ws.Range("U2").Value = 0.001879124199714181
$ws.Range("V2").Value = 0.001879124199714181
$ws.Range("W2").Value = 0.001834861700174168
$ws.Range("X2").Value = 0.001811552753176238
$ws.Range("Y2").Value = 0.001799634238166504

$ws.Range("C3").Value = 0.2319471836090088
$ws.Range("E3").Value = 94.10179021009753
$ws.Range("F3").Value = 0.003757337172971522
$ws.Range("G3").Value = 0.00300435408836291
$ws.Range("H3").Value = 0.002929127969900113
$ws.Range("I3").Value = 0.002475309503425576
$ws.Range("J3").Value = 0.002448835996969612
$ws.Range("K3").Value = 0.002210805743746385
$ws.Range("L3").Value = 0.002101490780628188
$ws.Range("M3").Value = 0.002101490780628188
$ws.Range("N3").Value = 0.00199639128112805
$ws.Range("O3").Value = 0.001980608108248335
$ws.Range("P3").Value = 0.001980608108248335
$ws.Range("Q3").Value = 0.001980608108248335
$ws.Range("R3").Value = 0.001980608108248335
$ws.Range("S3").Value = 0.001935560224642048
$ws.Range("T3").Value = 0.001892956354633646
$ws.Range("U3").Value = 0.001892956354633646
$ws.Range("V3").Value = 0.001877605169704809
$ws.Range("W3").Value = 0.001859615969279884
$ws.Range("X3").Value = 0.001834342889085722
$ws.Range("Y3").Value = 0.001834342889085722

$ws.Range("C4").Value = 0.2222318649291992
$ws.Range("E4").Value = 93.56065908049823
$ws.Range("F4").Value = 0.003468199983191212
$ws.Range("G4").Value = 0.002873950838637925
$ws.Range("H4").Value = 0.002703252195435417
$ws.Range("I4").Value = 0.00253146159266655
$ws.Range("J4").Value = 0.002453902374972207
$ws.Range("K4").Value = 0.002237395407583961
$ws.Range("L4").Value = 0.002237395407583961
$ws.Range("M4").Value = 0.002105108836197875
$ws.Range("N4").Value = 0.002062667900228055
$ws.Range("O4").Value = 0.002062667900228055
$ws.Range("P4").Value = 0.001960094144735697
$ws.Range("Q4").Value = 0.00195714878499223
$ws.Range("R4").Value = 0.001913328637333817
$ws.Range("S4").Value = 0.001910482265665615
$ws.Range("T4").Value = 0.001892134270515652
$ws.Range("U4").Value = 0.001892134270515652
$ws.Range("V4").Value = 0.001873776959767775
$ws.Range("W4").Value = 0.001845402509456848
$ws.Range("X4").Value = 0.001838388140349565
$ws.Range("Y4").Value = 0.00182379452398632

$ws.Range("C5").Value = 0.2114667892456055
$ws.Range("E5").Value = 92.85271622000255
$ws.Range("F5").Value = 0.003757337172971522
$ws.Range("G5").Value = 0.002907773767802515
$ws.Range("H5").Value = 0.00275645264018522
$ws.Range("I5").Value = 0.002519862808729168
$ws.Range("J5").Value = 0.002263607809034166
$ws.Range("K5").Value = 0.002207013663649165
$ws.Range("L5").Value = 0.002168499731804785
$ws.Range("M5").Value = 0.001995674515961284
$ws.Range("N5").Value = 0.001995674515961284
$ws.Range("O5").Value = 0.001995674515961284
$ws.Range("P5").Value = 0.001995674515961284
$ws.Range("Q5").Value = 0.001988859831540857
$ws.Range("R5").Value = 0.001950726211998344
$ws.Range("S5").Value = 0.001932168180179115
$ws.Range("T5").Value = 0.001916603307719151
$ws.Range("U5").Value = 0.001916603307719151
$ws.Range("V5").Value = 0.001886347806190123
$ws.Range("W5").Value = 0.001845638826943096
$ws.Range("X5").Value = 0.001810904081250296
$ws.Range("Y5").Value = 0.00180999446822617

$ws.Range("C6").Value = 0.2207469940185547
$ws.Range("E6").Value = 95.00142985280218
$ws.Range("F6").Value = 0.003626117340231359
$ws.Range("G6").Value = 0.002941968445532599
$ws.Range("H6").Value = 0.002480681334198425
$ws.Range("I6").Value = 0.002480681334198425
$ws.Range("J6").Value = 0.002382963304980716
$ws.Range("K6").Value = 0.002308582547353814
$ws.Range("L6").Value = 0.002242529045694349
$ws.Range("M6").Value = 0.002242529045694349
$ws.Range("N6").Value = 0.002122802250367839
$ws.Range("O6").Value = 0.002048204246051749
$ws.Range("P6").Value = 0.002040680870745346
$ws.Range("Q6").Value = 0.001997850410463999
$ws.Range("R6").Value = 0.001997850410463999
$ws.Range("S6").Value = 0.001925365160903658
$ws.Range("T6").Value = 0.001925365160903658
$ws.Range("U6").Value = 0.001925365160903658
$ws.Range("V6").Value = 0.001877031183387471
$ws.Range("W6").Value = 0.001866983752396668
$ws.Range("X6").Value = 0.001864732699654495
$ws.Range("Y6").Value = 0.001851879724226163

$ws.Range("C7").Value = 0.2109410762786865
$ws.Range("E7").Value = 97.39537783203195
$ws.Range("G7").Value = 0.002842662240405421
$ws.Range("H7").Value = 0.002590884455268118
$ws.Range("I7").Value = 0.002570799048397084
$ws.Range("J7").Value = 0.002421251962234598
$ws.Range("K7").Value = 0.002326692770400921
$ws.Range("L7").Value = 0.002073176551154023
$ws.Range("M7").Value = 0.002073176551154023
$ws.Range("N7").Value = 0.002073176551154023
$ws.Range("O7").Value = 0.002073176551154023
$ws.Range("P7").Value = 0.002004578295734327
$ws.Range("Q7").Value = 0.002004578295734327
$ws.Range("R7").Value = 0.002004578295734327
$ws.Range("S7").Value = 0.002004578295734327
$ws.Range("T7").Value = 0.001999382805630485
$ws.Range("U7").Value = 0.001999382805630485
$ws.Range("V7").Value = 0.00198616754964539
$ws.Range("W7").Value = 0.001938230577194991
$ws.Range("X7").Value = 0.001912721410545843
$ws.Range("Y7").Value = 0.001898545376842728

$ws.Range("C8").Value = 0.2369227409362793
$ws.Range("E8").Value = 94.43850981941432
$ws.Range("F8").Value = 0.003700162937976026
$ws.Range("G8").Value = 0.003009598797477282
$ws.Range("H8").Value = 0.002646269374357229
$ws.Range("I8").Value = 0.002552568549647183
$ws.Range("J8").Value = 0.00244199668874411
$ws.Range("K8").Value = 0.002381499876536922
$ws.Range("L8").Value = 0.002381499876536922
$ws.Range("M8").Value = 0.002348030297378829
$ws.Range("N8").Value = 0.002281638363969264
$ws.Range("O8").Value = 0.002263781616278699
$ws.Range("P8").Value = 0.002238783968251937
$ws.Range("Q8").Value = 0.002116550579053741
$ws.Range("R8").Value = 0.002006542688852196
$ws.Range("S8").Value = 0.001950845564374991
$ws.Range("T8").Value = 0.00194626803090169
$ws.Range("U8").Value = 0.001909999530603192
$ws.Range("V8").Value = 0.001857911206957919
$ws.Range("W8").Value = 0.001857911206957919
$ws.Range("X8").Value = 0.001847654058312319
$ws.Range("Y8").Value = 0.001840906624160123

$ws.Range("C9").Value = 0.2717306613922119
$ws.Range("E9").Value = 91.59437300454374
$ws.Range("F9").Value = 0.003757337172971522
$ws.Range("G9").Value = 0.003001218712700802
$ws.Range("H9").Value = 0.002907999265922204
$ws.Range("I9").Value = 0.002496850785575897
$ws.Range("J9").Value = 0.002313348660119449
$ws.Range("K9").Value = 0.002313348660119449
$ws.Range("L9").Value = 0.002313348660119449
$ws.Range("M9").Value = 0.002198293861550369
$ws.Range("N9").Value = 0.002084392455973664
$ws.Range("O9").Value = 0.002042702684759833
$ws.Range("P9").Value = 0.001940469004323559
$ws.Range("Q9").Value = 0.001894629648086088
$ws.Range("R9").Value = 0.001894629648086088
$ws.Range("S9").Value = 0.001867621639977022
$ws.Range("T9").Value = 0.001824675259167616
$ws.Range("U9").Value = 0.001824675259167616
$ws.Range("V9").Value = 0.001824675259167616
$ws.Range("W9").Value = 0.001802788024790352
$ws.Range("X9").Value = 0.001787784331931919
$ws.Range("Y9").Value = 0.001785465360712353

$ws.Range("C10").Value = 0.3041648864746094
$ws.Range("E10").Value = 97.24944864968529
$ws.Range("F10").Value = 0.003660312753857895
$ws.Range("G10").Value = 0.003041080166472455
$ws.Range("H10").Value = 0.002764988425102703
$ws.Range("I10").Value = 0.002383870006744072
$ws.Range("J10").Value = 0.002383870006744072
$ws.Range("K10").Value = 0.002139053313738748
$ws.Range("L10").Value = 0.002139053313738748
$ws.Range("M10").Value = 0.002118432442789481
$ws.Range("N10").Value = 0.002039869460326829
$ws.Range("O10").Value = 0.002039869460326829
$ws.Range("P10").Value = 0.002039869460326829
$ws.Range("Q10").Value = 0.002022814060199
$ws.Range("R10").Value = 0.001971363524149658
$ws.Range("S10").Value = 0.001971363524149658
$ws.Range("T10").Value = 0.001971363524149658
$ws.Range("U10").Value = 0.001971363524149658
$ws.Range("V10").Value = 0.001945630741105822
$ws.Range("W10").Value = 0.001908769368766341
$ws.Range("X10").Value = 0.001908769368766341
$ws.Range("Y10").Value = 0.001895700753405171

$ws.Range("C11").Value = 0.2020366191864014
$ws.Range("E11").Value = 90.00817902700874
$ws.Range("F11").Value = 0.003669970315213592
$ws.Range("G11").Value = 0.002708987109506967
$ws.Range("H11").Value = 0.00252902655463289
$ws.Range("I11").Value = 0.002462306971627018
$ws.Range("J11").Value = 0.002236275627551824
$ws.Range("K11").Value = 0.002179682682382871
$ws.Range("L11").Value = 0.002078025384390665
$ws.Range("M11").Value = 0.001946477616925514
$ws.Range("N11").Value = 0.001946477616925514
$ws.Range("O11").Value = 0.001946477616925514
$ws.Range("P11").Value = 0.001946477616925514
$ws.Range("Q11").Value = 0.001946477616925514
$ws.Range("R11").Value = 0.001803422409534726
$ws.Range("S11").Value = 0.001803422409534726
$ws.Range("T11").Value = 0.001803422409534726
$ws.Range("U11").Value = 0.001803422409534726
$ws.Range("V11").Value = 0.001803422409534726
$ws.Range("W11").Value = 0.001758293220352383
$ws.Range("X11").Value = 0.001758293220352383
$ws.Range("Y11").Value = 0.001754545400136622
